# Backup QR Scanner data refresh (2025-10-02T16:26:33.634Z)
# Replaces the Chest-scanner log rows with a fresh pull: updated Student ID
# order, Log Date rolled from 01/10/2025 to 02/10/2025, new Log Time stamps,
# and 9 additional rows (40-48) appended from the new scan session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Student ID / Log Date / Log Time are stored as text (not numbers/dates),
# so force the "@" text format before writing values into these columns.
$ws.Range("A2:A48").NumberFormat = "@"
$ws.Range("C2:C48").NumberFormat = "@"
$ws.Range("D2:D48").NumberFormat = "@"
$ws.Range("F40:F48").NumberFormat = "@"

$ws.Range("A2").Value = "221672"
$ws.Range("C2").Value = "02/10/2025"
$ws.Range("D2").Value = "19:17:37"
$ws.Range("A3").Value = "221713"
$ws.Range("C3").Value = "02/10/2025"
$ws.Range("D3").Value = "19:17:44"
$ws.Range("A4").Value = "221582"
$ws.Range("C4").Value = "02/10/2025"
$ws.Range("D4").Value = "19:17:56"
$ws.Range("A5").Value = "221535"
$ws.Range("C5").Value = "02/10/2025"
$ws.Range("D5").Value = "19:18:05"
$ws.Range("A6").Value = "221676"
$ws.Range("C6").Value = "02/10/2025"
$ws.Range("D6").Value = "19:18:12"
$ws.Range("A7").Value = "221722"
$ws.Range("C7").Value = "02/10/2025"
$ws.Range("D7").Value = "19:18:20"
$ws.Range("A8").Value = "221712"
$ws.Range("C8").Value = "02/10/2025"
$ws.Range("D8").Value = "19:18:49"
$ws.Range("A9").Value = "221641"
$ws.Range("C9").Value = "02/10/2025"
$ws.Range("D9").Value = "19:19:01"
$ws.Range("A10").Value = "221655"
$ws.Range("C10").Value = "02/10/2025"
$ws.Range("D10").Value = "19:19:10"
$ws.Range("A11").Value = "221605"
$ws.Range("C11").Value = "02/10/2025"
$ws.Range("D11").Value = "19:19:21"
$ws.Range("A12").Value = "221592"
$ws.Range("C12").Value = "02/10/2025"
$ws.Range("D12").Value = "19:19:29"
$ws.Range("A13").Value = "221658"
$ws.Range("C13").Value = "02/10/2025"
$ws.Range("D13").Value = "19:19:40"
$ws.Range("A14").Value = "221552"
$ws.Range("C14").Value = "02/10/2025"
$ws.Range("D14").Value = "19:19:54"
$ws.Range("A15").Value = "221752"
$ws.Range("C15").Value = "02/10/2025"
$ws.Range("D15").Value = "19:20:10"
$ws.Range("A16").Value = "221751"
$ws.Range("C16").Value = "02/10/2025"
$ws.Range("D16").Value = "19:20:20"
$ws.Range("A17").Value = "221706"
$ws.Range("C17").Value = "02/10/2025"
$ws.Range("D17").Value = "19:20:30"
$ws.Range("A18").Value = "221568"
$ws.Range("C18").Value = "02/10/2025"
$ws.Range("D18").Value = "19:20:40"
$ws.Range("A19").Value = "221719"
$ws.Range("C19").Value = "02/10/2025"
$ws.Range("D19").Value = "19:21:08"
$ws.Range("A20").Value = "221657"
$ws.Range("C20").Value = "02/10/2025"
$ws.Range("D20").Value = "19:21:23"
$ws.Range("A21").Value = "221596"
$ws.Range("C21").Value = "02/10/2025"
$ws.Range("D21").Value = "19:21:36"
$ws.Range("A22").Value = "221621"
$ws.Range("C22").Value = "02/10/2025"
$ws.Range("D22").Value = "19:21:50"
$ws.Range("A23").Value = "221538"
$ws.Range("C23").Value = "02/10/2025"
$ws.Range("D23").Value = "19:22:03"
$ws.Range("A24").Value = "221595"
$ws.Range("C24").Value = "02/10/2025"
$ws.Range("D24").Value = "19:22:11"
$ws.Range("A25").Value = "221642"
$ws.Range("C25").Value = "02/10/2025"
$ws.Range("D25").Value = "19:22:18"
$ws.Range("A26").Value = "221556"
$ws.Range("C26").Value = "02/10/2025"
$ws.Range("D26").Value = "19:22:44"
$ws.Range("A27").Value = "221624"
$ws.Range("C27").Value = "02/10/2025"
$ws.Range("D27").Value = "19:22:51"
$ws.Range("A28").Value = "221745"
$ws.Range("C28").Value = "02/10/2025"
$ws.Range("D28").Value = "19:23:14"
$ws.Range("A29").Value = "221594"
$ws.Range("C29").Value = "02/10/2025"
$ws.Range("D29").Value = "19:23:23"
$ws.Range("A30").Value = "221547"
$ws.Range("C30").Value = "02/10/2025"
$ws.Range("D30").Value = "19:23:30"
$ws.Range("A31").Value = "221533"
$ws.Range("C31").Value = "02/10/2025"
$ws.Range("D31").Value = "19:23:41"
$ws.Range("A32").Value = "221674"
$ws.Range("C32").Value = "02/10/2025"
$ws.Range("D32").Value = "19:23:51"
$ws.Range("A33").Value = "221697"
$ws.Range("C33").Value = "02/10/2025"
$ws.Range("D33").Value = "19:24:00"
$ws.Range("A34").Value = "221683"
$ws.Range("C34").Value = "02/10/2025"
$ws.Range("D34").Value = "19:24:08"
$ws.Range("A35").Value = "221694"
$ws.Range("C35").Value = "02/10/2025"
$ws.Range("D35").Value = "19:24:15"
$ws.Range("A36").Value = "221709"
$ws.Range("C36").Value = "02/10/2025"
$ws.Range("D36").Value = "19:24:28"
$ws.Range("A37").Value = "221707"
$ws.Range("C37").Value = "02/10/2025"
$ws.Range("D37").Value = "19:24:36"
$ws.Range("A38").Value = "221689"
$ws.Range("C38").Value = "02/10/2025"
$ws.Range("D38").Value = "19:24:43"
$ws.Range("A39").Value = "221716"
$ws.Range("C39").Value = "02/10/2025"
$ws.Range("D39").Value = "19:24:49"
$ws.Range("A40").Value = "221686"
$ws.Range("B40").Value = "Chest"
$ws.Range("C40").Value = "02/10/2025"
$ws.Range("D40").Value = "19:25:05"
$ws.Range("E40").Value = "Manual"
$ws.Range("F40").Value = "ahmedali78112@gmail.com"
$ws.Range("A41").Value = "221711"
$ws.Range("B41").Value = "Chest"
$ws.Range("C41").Value = "02/10/2025"
$ws.Range("D41").Value = "19:25:12"
$ws.Range("E41").Value = "Manual"
$ws.Range("F41").Value = "ahmedali78112@gmail.com"
$ws.Range("A42").Value = "221673"
$ws.Range("B42").Value = "Chest"
$ws.Range("C42").Value = "02/10/2025"
$ws.Range("D42").Value = "19:25:19"
$ws.Range("E42").Value = "Manual"
$ws.Range("F42").Value = "ahmedali78112@gmail.com"
$ws.Range("A43").Value = "221616"
$ws.Range("B43").Value = "Chest"
$ws.Range("C43").Value = "02/10/2025"
$ws.Range("D43").Value = "19:25:40"
$ws.Range("E43").Value = "Manual"
$ws.Range("F43").Value = "ahmedali78112@gmail.com"
$ws.Range("A44").Value = "221542"
$ws.Range("B44").Value = "Chest"
$ws.Range("C44").Value = "02/10/2025"
$ws.Range("D44").Value = "19:25:48"
$ws.Range("E44").Value = "Manual"
$ws.Range("F44").Value = "ahmedali78112@gmail.com"
$ws.Range("A45").Value = "221585"
$ws.Range("B45").Value = "Chest"
$ws.Range("C45").Value = "02/10/2025"
$ws.Range("D45").Value = "19:25:57"
$ws.Range("E45").Value = "Manual"
$ws.Range("F45").Value = "ahmedali78112@gmail.com"
$ws.Range("A46").Value = "221603"
$ws.Range("B46").Value = "Chest"
$ws.Range("C46").Value = "02/10/2025"
$ws.Range("D46").Value = "19:26:05"
$ws.Range("E46").Value = "Manual"
$ws.Range("F46").Value = "ahmedali78112@gmail.com"
$ws.Range("A47").Value = "221677"
$ws.Range("B47").Value = "Chest"
$ws.Range("C47").Value = "02/10/2025"
$ws.Range("D47").Value = "19:26:11"
$ws.Range("E47").Value = "Manual"
$ws.Range("F47").Value = "ahmedali78112@gmail.com"
$ws.Range("A48").Value = "180915"
$ws.Range("B48").Value = "Chest"
$ws.Range("C48").Value = "02/10/2025"
$ws.Range("D48").Value = "19:26:18"
$ws.Range("E48").Value = "Manual"
$ws.Range("F48").Value = "ahmedali78112@gmail.com"
